$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price cells keep their literal text representation
# (values like "0.5220", "307.05", "27.134.55" must stay text, not become numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.134.55"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.899.88"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.05"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5220"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3803"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07287"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.34"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9024"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08169"
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.35"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.851.53"
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.352"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008647"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.67"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.175.56"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.118"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.80"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.457"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.329"
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.11"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.25"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.735"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.65"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.894"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09211"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05039"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7938"
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.218"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.970"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.363"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.631"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5710"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.080"
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.044"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.598"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.59"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1512"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4882"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.19"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.632"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.38"
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.90"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("E51").Value = "  +0.55%  "
